# Add "wants" (Not-evening / Not-task constraint lists) for operators in the
# Operators sheet. This populates two new columns (G: "Not evening",
# H: "Not task") with per-operator comma-separated lists of day/task numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Operators")
$ws.Activate() | Out-Null

# Header row
$ws.Cells.Item(1, 7).Value = "Not evening"
$ws.Cells.Item(1, 8).Value = "Not task"

# Row 2 (O1)
$ws.Cells.Item(2, 8).Value = "26,27"

# Row 3 (O2)
$ws.Cells.Item(3, 8).Value = "8,12,13,26,27"

# Row 4 (O3)
$ws.Cells.Item(4, 7).Value = "2,9,18,23,30"
$ws.Cells.Item(4, 8).Value = "5,6,19,20"

# Row 5 (O4)
$ws.Cells.Item(5, 7).Value = "4,9,16"
$ws.Cells.Item(5, 8).Value = "5,6,26,27"

# Row 6 (O5)
$ws.Cells.Item(6, 7).Value = "9,25"
$ws.Cells.Item(6, 8).Value = "26,27"

# Row 7 (O6)
$ws.Cells.Item(7, 7).Value = "9,18,23,30"
$ws.Cells.Item(7, 8).Value = "5,6,19,20"

# Row 8 (O7)
$ws.Cells.Item(8, 7).Value = "2,9,15,23"
$ws.Cells.Item(8, 8).Value = "26,27"

# Row 9 (O8)
$ws.Cells.Item(9, 7).Value = "2,9"
$ws.Cells.Item(9, 8).Value = "5,6,12,13"

# Row 10 (O9)
$ws.Cells.Item(10, 7).Value = "4,9,15,22,29"
$ws.Cells.Item(10, 8).Value = "5,6,26,27"

# Row 11 (O10)
$ws.Cells.Item(11, 7).Value = "2,11,17,23"
$ws.Cells.Item(11, 8).Value = "1,3,12,13,24"

# Row 12 (O11)
$ws.Cells.Item(12, 7).Value = 23
$ws.Cells.Item(12, 8).Value = "12,13,26,27"

# Row 13 (O12)
$ws.Cells.Item(13, 7).Value = "11,14,23"
$ws.Cells.Item(13, 8).Value = "5,6,26,27"

# Row 14 (O13)
$ws.Cells.Item(14, 7).Value = 23
$ws.Cells.Item(14, 8).Value = "26,27,28"

# Row 15 (O14)
$ws.Cells.Item(15, 8).Value = "11,12,13"

# Row 16 (O15)
$ws.Cells.Item(16, 7).Value = "14,25,31"
$ws.Cells.Item(16, 8).Value = "9,12,13,26,27"

# Row 17 (O16)
$ws.Cells.Item(17, 8).Value = "2,3,4,5,6,8,10,11,17,23,25,26,27"

# Row 18 (O17)
$ws.Cells.Item(18, 7).Value = 8
$ws.Cells.Item(18, 8).Value = "9,10,11,14,15,16,17,18,19,20,21,23,25,26,27,28"

# Row 19 (O18)
$ws.Cells.Item(19, 7).Value = "2,9"
$ws.Cells.Item(19, 8).Value = "11,12,13,14,26,27"

# Row 20 (O19)
$ws.Cells.Item(20, 7).Value = "3,23,29"
$ws.Cells.Item(20, 8).Value = "7,17,26,27"

# Row 22 (O21)
$ws.Cells.Item(22, 7).Value = 10
$ws.Cells.Item(22, 8).Value = "9,17,26,27"

# Row 23 (O22)
$ws.Cells.Item(23, 7).Value = "10,18,23"
$ws.Cells.Item(23, 8).Value = "5,6,19,20"

# Row 24 (O23)
$ws.Cells.Item(24, 8).Value = "2,3,5,6,10,21,26,27"

# Row 25 (O24)
$ws.Cells.Item(25, 7).Value = "17,25"
$ws.Cells.Item(25, 8).Value = "3,10,21,23,26,27,28"

# Row 26 (O25)
$ws.Cells.Item(26, 8).Value = "8,10,11,12,13,25,26,27"

# Row 27 (O26)
$ws.Cells.Item(27, 8).Value = "19,20,26,27,30"

# Row 28 (O27)
$ws.Cells.Item(28, 7).Value = 10
$ws.Cells.Item(28, 8).Value = "19,20,26,27"

# Row 29 (O28)
$ws.Cells.Item(29, 7).Value = 7
$ws.Cells.Item(29, 8).Value = "3,5,6,17,25,26,27"

# Row 30 (O29)
$ws.Cells.Item(30, 7).Value = "7,15,23"
$ws.Cells.Item(30, 8).Value = "10,11,26,27,28"

# Row 31 (O30)
$ws.Cells.Item(31, 7).Value = 7
$ws.Cells.Item(31, 8).Value = "10,11,18,19,20,21,26,27"

# Row 32 (O31)
$ws.Cells.Item(32, 7).Value = "7,23"
$ws.Cells.Item(32, 8).Value = "10,11,26,27,28"

# Row 33 (O32)
$ws.Cells.Item(33, 8).Value = "1,26,27"

# Row 34 (O33)
$ws.Cells.Item(34, 7).Value = 7
$ws.Cells.Item(34, 8).Value = "3,14,17"

# Row 35 (O34)
$ws.Cells.Item(35, 8).Value = "1,3,4,5,6,9,16,23,26,27"

# Row 38 (O37)
$ws.Cells.Item(38, 7).Value = 25
$ws.Cells.Item(38, 8).Value = "9,19,20,26,27,28"

# Row 40 (O39)
$ws.Cells.Item(40, 8).Value = "26,27,28"

# Row 41 (O40)
$ws.Cells.Item(41, 7).Value = 18

# Row 42 (O41)
$ws.Cells.Item(42, 8).Value = "1,11"

# Row 43 (O42)
$ws.Cells.Item(43, 8).Value = "2,9,16,23"

# New column widths to mirror the authored layout (bestFit columns G/H)
$ws.Columns.Item(7).ColumnWidth = 9.36
$ws.Columns.Item(8).ColumnWidth = 6.36

# Match the author's final selection/view state (Operators tab active, G12 selected)
$ws.Range("G12").Select() | Out-Null
